$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column H (8th column) to match the target layout
$ws.Columns(8).ColumnWidth = 26.285714285714285

# New "data_type" column header, styled like the other header cells
$ws.Range("I1").Value = "data_type"
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)

# Per-row data type classification for column I
$types = @("int","float","int","float","float","bool","bool","bool","bool","float","float","bool","float","bool","float","float","int","int","int","int","int","int","int","int","int","float","float","int","int","int","int","int","bool","bool","float","float","bool","int","float","float","string","string","int","bool","bool","bool","bool","bool","bool","bool","int","int","int","string","string","string","string","string","string","int","int","string","string","bool","string","bool","bool","bool","string","string","string","string","string","string","bool","bool","string","bool","bool","string","bool","bool","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","bool","bool","string","string","string","string","string","string","string","string","string","string","string","string","int","string","string","string","string","string","string","int","int","string","string","bool","string","bool","bool","bool","string","string","string","string","string","string","bool","bool","string","bool","bool","string","bool","bool","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","string","float","float","float","float","float","float","float","float","float","float","int","int","float","float","float","float","float","float","float")
for ($i = 0; $i -lt $types.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $types[$i]
}

# D159: survey question type corrected from Multiple_choice to Scale
$ws.Range("D159").Value = "Scale"

# Update view state to match the saved selection/scroll position
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("I2").Select()
